# Update countries & provincias Spain
#
# - Refresh the "last updated" timestamp (17:04 -> 18:21)
# - Refresh the COVID counters for a batch of countries (new totals,
#   new cases, active cases, recovered, deaths today, deaths)
# - Two pairs/triples of countries had swapped positions in the sheet:
#     * Grecia / Finlandia (rows 99/100) - Grecia now comes first with the
#       refreshed figures, Finlandia drops to row 100 keeping its old figures
#     * Lesoto / Reunion / Republica del Chad (rows 159/160/161) - Lesoto and
#       Reunion move up (with refreshed figures), Republica del Chad moves
#       down keeping its old figures

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner
$ws.Range("A1").Value = 'Datos actualizados a 20 de Agosto de 2020 a las 18:21'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5713181
$ws.Range("C4").Value = 12250
$ws.Range("D4").Value = 3064213
$ws.Range("E4").Value = 2472303
$ws.Range("G4").Value = 328
$ws.Range("H4").Value = 176665

# India (row 6)
$ws.Range("B6").Value = 2900967
$ws.Range("C6").Value = 65145
$ws.Range("D6").Value = 2153342
$ws.Range("E6").Value = 692677
$ws.Range("G6").Value = 954
$ws.Range("H6").Value = 54948

# Chile (row 12)
$ws.Range("B12").Value = 391849
$ws.Range("C12").Value = 1812
$ws.Range("E12").Value = 16893
$ws.Range("G12").Value = 93
$ws.Range("H12").Value = 10671

# Reino Unido (row 15)
$ws.Range("B15").Value = 322280
$ws.Range("C15").Value = 1182
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 41403

# Italia (row 20)
$ws.Range("B20").Value = 256118
$ws.Range("C20").Value = 840
$ws.Range("D20").Value = 204686
$ws.Range("E20").Value = 16014
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 35418

# Alemania (row 22)
$ws.Range("B22").Value = 230018
$ws.Range("C22").Value = 318
$ws.Range("D22").Value = 204800
$ws.Range("E22").Value = 15902
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 9316

# Canada (row 27)
$ws.Range("B27").Value = 123653
$ws.Range("C27").Value = 163
$ws.Range("D27").Value = 110049
$ws.Range("E27").Value = 4553
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 9051

# Chequia (row 74)
$ws.Range("B74").Value = 20967
$ws.Range("C74").Value = 169
$ws.Range("D74").Value = 15852
$ws.Range("E74").Value = 4710
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 405

# Camerun (row 75)
$ws.Range("B75").Value = 18762
$ws.Range("C75").Value = 138
$ws.Range("E75").Value = 1814
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 408

# Noruega (row 89)
$ws.Range("D89").Value = 9150
$ws.Range("E89").Value = 776

# Grecia / Finlandia swap (rows 99 & 100)
$ws.Range("A99").Value = 'Grecia'
$ws.Range("B99").Value = 7934
$ws.Range("C99").Value = 250
$ws.Range("D99").Value = 3804
$ws.Range("E99").Value = 3895
$ws.Range("H99").Value = 235

$ws.Range("A100").Value = 'Finlandia'
$ws.Range("B100").Value = 7842
$ws.Range("C100").Value = 37
$ws.Range("D100").Value = 7100
$ws.Range("E100").Value = 408
$ws.Range("H100").Value = 334

# Sri Lanka (row 125)
$ws.Range("B125").Value = 2918
$ws.Range("C125").Value = 16
$ws.Range("E125").Value = 142

# Liberia (row 153)
$ws.Range("B153").Value = 1284
$ws.Range("C153").Value = 2
$ws.Range("E153").Value = 399

# Lesoto / Reunion / Republica del Chad rotation (rows 159, 160 & 161)
$ws.Range("A159").Value = 'Lesoto'
$ws.Range("B159").Value = 996
$ws.Range("C159").Value = 50
$ws.Range("D159").Value = 423
$ws.Range("E159").Value = 543
$ws.Range("H159").Value = 30

$ws.Range("A160").Value = 'Reunion'
$ws.Range("B160").Value = 996
$ws.Range("C160").Value = 51
$ws.Range("D160").Value = 657
$ws.Range("E160").Value = 334
$ws.Range("H160").Value = 5

$ws.Range("A161").Value = 'Republica del Chad'
$ws.Range("B161").Value = 971
$ws.Range("D161").Value = 868
$ws.Range("E161").Value = 27
$ws.Range("H161").Value = 76

# Belice (row 167)
$ws.Range("E167").Value = 510
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 5

# Comoras (row 172)
$ws.Range("B172").Value = 417
$ws.Range("C172").Value = 11
$ws.Range("D172").Value = 396
$ws.Range("E172").Value = 14
